# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, ± margins, dollar amounts) inside specific bullet paragraphs
# of the resume, matching the author's commit:
#   "Implement quantitative metrics highlighting across all resume formats"
#
# Approach: for each target paragraph, locate the metric substrings (in
# left-to-right order so repeated values resolve to the correct occurrence),
# compute absolute document character offsets from the paragraph's own
# Range.Start, and apply Bold + the "2C3E50" font color to just that
# sub-range via $d.Range(start, end). Word automatically splits/creates the
# runs needed (and adds xml:space="preserve" where appropriate), mirroring
# exactly what the diff shows.

$d = $word.ActiveDocument

# Hex 2C3E50 expressed the way Word's Font.Color (a COM RGB() value) expects
# it: RGB(0x2C, 0x3E, 0x50) = 0x2C + 0x3E*256 + 0x50*65536
$HighlightColor = 5258796

function Set-MetricHighlights {
    param($doc, $paraIndex, [string[]]$targets)

    $para = $doc.Paragraphs.Item($paraIndex)
    $range = $para.Range
    $fullText = $range.Text
    $searchFrom = 0

    foreach ($target in $targets) {
        $idx = $fullText.IndexOf($target, $searchFrom)
        if ($idx -lt 0) {
            $idx = $fullText.IndexOf($target)
        }
        if ($idx -ge 0) {
            $start = $range.Start + $idx
            $end = $start + $target.Length
            $sub = $doc.Range($start, $end)
            $sub.Font.Bold = 1
            $sub.Font.Color = $HighlightColor
            $searchFrom = $idx + $target.Length
        }
    }
}

# Paragraph 10: "• Discovered systematic race coding errors ... from 23% to 64%"
Set-MetricHighlights $d 10 @("23%", "64%")

# Paragraph 12: "• Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
Set-MetricHighlights $d 12 @("±4.2%", "±2.1%", "71%", "87%")

# Paragraph 13: "• Trigonometric algorithm ... 73.5% ... $4.7M ..."
Set-MetricHighlights $d 13 @("73.5%", "$4.7M")

# Paragraph 14: "• Built real-time FEC analysis systems ... valued over $2 trillion"
Set-MetricHighlights $d 14 @("$2")

# Paragraph 50: "• Predictive excellence: ... ±4.2% to ±2.1%"
Set-MetricHighlights $d 50 @("±4.2%", "±2.1%")

# Paragraph 51: "• Increased voter turnout prediction accuracy from 71% to 87%"
Set-MetricHighlights $d 51 @("71%", "87%")

# Paragraph 53: "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"
Set-MetricHighlights $d 53 @("34%", "28%")

Write-Output "Metrics highlighting applied"
